$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range('A1').Value = 'Datos actualizados a 21 de Marzo de 2020 a las 20:46'

# Country data table (rows 4-190), sorted descending by "Casos totales" (col B),
# reflecting the refreshed figures and resulting re-sort of the data.
$data = New-Object "object[,]" 187,8
$data[0,0] = 'China'
$data[0,1] = 81008
$data[0,2] = 41
$data[0,3] = 71740
$data[0,4] = 6013
$data[0,5] = 1927
$data[0,6] = 7
$data[0,7] = 3255
$data[1,0] = 'Italia'
$data[1,1] = 53578
$data[1,2] = 6557
$data[1,3] = 6072
$data[1,4] = 42681
$data[1,5] = 2857
$data[1,6] = 793
$data[1,7] = 4825
$data[2,0] = 'España'
$data[2,1] = 25374
$data[2,2] = 3803
$data[2,3] = 2125
$data[2,4] = 21871
$data[2,5] = 1612
$data[2,6] = 285
$data[2,7] = 1378
$data[3,0] = 'Estados Unidos'
$data[3,1] = 24137
$data[3,2] = 4754
$data[3,3] = 171
$data[3,4] = 23678
$data[3,5] = 64
$data[3,6] = 32
$data[3,7] = 288
$data[4,0] = 'Alemania'
$data[4,1] = 22084
$data[4,2] = 2236
$data[4,3] = 209
$data[4,4] = 21792
$data[4,5] = 2
$data[4,6] = 15
$data[4,7] = 83
$data[5,0] = 'Iran'
$data[5,1] = 20610
$data[5,2] = 966
$data[5,3] = 7635
$data[5,4] = 11419
$data[5,5] = 0
$data[5,6] = 123
$data[5,7] = 1556
$data[6,0] = 'Francia'
$data[6,1] = 14459
$data[6,2] = 1847
$data[6,3] = 1587
$data[6,4] = 12310
$data[6,5] = 1525
$data[6,6] = 112
$data[6,7] = 562
$data[7,0] = 'Corea del Sur'
$data[7,1] = 8799
$data[7,2] = 147
$data[7,3] = 2612
$data[7,4] = 6085
$data[7,5] = 59
$data[7,6] = 8
$data[7,7] = 102
$data[8,0] = 'Suiza'
$data[8,1] = 6665
$data[8,2] = 1050
$data[8,3] = 15
$data[8,4] = 6575
$data[8,5] = 141
$data[8,6] = 19
$data[8,7] = 75
$data[9,0] = 'Reino Unido'
$data[9,1] = 5018
$data[9,2] = 1035
$data[9,3] = 65
$data[9,4] = 4720
$data[9,5] = 20
$data[9,6] = 56
$data[9,7] = 233
$data[10,0] = 'Paises Bajos'
$data[10,1] = 3631
$data[10,2] = 637
$data[10,3] = 2
$data[10,4] = 3493
$data[10,5] = 354
$data[10,6] = 30
$data[10,7] = 136
$data[11,0] = 'Austria'
$data[11,1] = 2847
$data[11,2] = 198
$data[11,3] = 9
$data[11,4] = 2830
$data[11,5] = 15
$data[11,6] = 2
$data[11,7] = 8
$data[12,0] = 'Belgica'
$data[12,1] = 2815
$data[12,2] = 558
$data[12,3] = 263
$data[12,4] = 2485
$data[12,5] = 288
$data[12,6] = 30
$data[12,7] = 67
$data[13,0] = 'Noruega'
$data[13,1] = 2141
$data[13,2] = 182
$data[13,3] = 1
$data[13,4] = 2133
$data[13,5] = 28
$data[13,6] = 0
$data[13,7] = 7
$data[14,0] = 'Suecia'
$data[14,1] = 1764
$data[14,2] = 125
$data[14,3] = 16
$data[14,4] = 1728
$data[14,5] = 69
$data[14,6] = 4
$data[14,7] = 20
$data[15,0] = 'Dinamarca'
$data[15,1] = 1326
$data[15,2] = 71
$data[15,3] = 1
$data[15,4] = 1312
$data[15,5] = 42
$data[15,6] = 4
$data[15,7] = 13
$data[16,0] = 'Portugal'
$data[16,1] = 1280
$data[16,2] = 260
$data[16,3] = 5
$data[16,4] = 1263
$data[16,5] = 26
$data[16,6] = 6
$data[16,7] = 12
$data[17,0] = 'Canada'
$data[17,1] = 1205
$data[17,2] = 118
$data[17,3] = 14
$data[17,4] = 1173
$data[17,5] = 1
$data[17,6] = 6
$data[17,7] = 18
$data[18,0] = 'Malasia'
$data[18,1] = 1183
$data[18,2] = 153
$data[18,3] = 114
$data[18,4] = 1061
$data[18,5] = 26
$data[18,6] = 5
$data[18,7] = 8
$data[19,0] = 'Australia'
$data[19,1] = 1072
$data[19,2] = 144
$data[19,3] = 46
$data[19,4] = 1019
$data[19,5] = 2
$data[19,6] = 0
$data[19,7] = 7
$data[20,0] = 'Japon'
$data[20,1] = 1046
$data[20,2] = 39
$data[20,3] = 215
$data[20,4] = 795
$data[20,5] = 55
$data[20,6] = 1
$data[20,7] = 36
$data[21,0] = 'Brasil'
$data[21,1] = 1021
$data[21,2] = 51
$data[21,3] = 2
$data[21,4] = 1001
$data[21,5] = 18
$data[21,6] = 7
$data[21,7] = 18
$data[22,0] = 'Chequia'
$data[22,1] = 995
$data[22,2] = 162
$data[22,3] = 6
$data[22,4] = 989
$data[22,5] = 7
$data[22,6] = 0
$data[22,7] = 0
$data[23,0] = 'Israel'
$data[23,1] = 883
$data[23,2] = 178
$data[23,3] = 36
$data[23,4] = 846
$data[23,5] = 15
$data[23,6] = 0
$data[23,7] = 1
$data[24,0] = 'Irlanda'
$data[24,1] = 785
$data[24,2] = 102
$data[24,3] = 5
$data[24,4] = 777
$data[24,5] = 6
$data[24,6] = 0
$data[24,7] = 3
$data[25,0] = 'Pakistan'
$data[25,1] = 734
$data[25,2] = 233
$data[25,3] = 13
$data[25,4] = 718
$data[25,5] = 0
$data[25,6] = 0
$data[25,7] = 3
$data[26,0] = 'Crucero'
$data[26,1] = 712
$data[26,2] = 0
$data[26,3] = 567
$data[26,4] = 137
$data[26,5] = 15
$data[26,6] = 0
$data[26,7] = 8
$data[27,0] = 'Turquia'
$data[27,1] = 670
$data[27,2] = 0
$data[27,3] = 0
$data[27,4] = 661
$data[27,5] = 0
$data[27,6] = 0
$data[27,7] = 9
$data[28,0] = 'Luxemburgo'
$data[28,1] = 670
$data[28,2] = 186
$data[28,3] = 6
$data[28,4] = 656
$data[28,5] = 3
$data[28,6] = 3
$data[28,7] = 8
$data[29,0] = 'Chile'
$data[29,1] = 537
$data[29,2] = 103
$data[29,3] = 8
$data[29,4] = 529
$data[29,5] = 7
$data[29,6] = 0
$data[29,7] = 0
$data[30,0] = 'Grecia'
$data[30,1] = 530
$data[30,2] = 35
$data[30,3] = 19
$data[30,4] = 498
$data[30,5] = 20
$data[30,6] = 3
$data[30,7] = 13
$data[31,0] = 'Finlandia'
$data[31,1] = 523
$data[31,2] = 73
$data[31,3] = 10
$data[31,4] = 512
$data[31,5] = 2
$data[31,6] = 1
$data[31,7] = 1
$data[32,0] = 'Ecuador'
$data[32,1] = 506
$data[32,2] = 80
$data[32,3] = 3
$data[32,4] = 496
$data[32,5] = 2
$data[32,6] = 0
$data[32,7] = 7
$data[33,0] = 'Polonia'
$data[33,1] = 492
$data[33,2] = 67
$data[33,3] = 13
$data[33,4] = 474
$data[33,5] = 3
$data[33,6] = 0
$data[33,7] = 5
$data[34,0] = 'Catar'
$data[34,1] = 481
$data[34,2] = 11
$data[34,3] = 27
$data[34,4] = 454
$data[34,5] = 6
$data[34,6] = 0
$data[34,7] = 0
$data[35,0] = 'Islandia'
$data[35,1] = 473
$data[35,2] = 64
$data[35,3] = 5
$data[35,4] = 467
$data[35,5] = 1
$data[35,6] = 1
$data[35,7] = 1
$data[36,0] = 'Indonesia'
$data[36,1] = 450
$data[36,2] = 81
$data[36,3] = 20
$data[36,4] = 392
$data[36,5] = 0
$data[36,6] = 6
$data[36,7] = 38
$data[37,0] = 'Singapur'
$data[37,1] = 432
$data[37,2] = 47
$data[37,3] = 140
$data[37,4] = 290
$data[37,5] = 14
$data[37,6] = 2
$data[37,7] = 2
$data[38,0] = 'Tailandia'
$data[38,1] = 411
$data[38,2] = 89
$data[38,3] = 44
$data[38,4] = 366
$data[38,5] = 7
$data[38,6] = 0
$data[38,7] = 1
$data[39,0] = 'Arabia Saudita'
$data[39,1] = 392
$data[39,2] = 48
$data[39,3] = 16
$data[39,4] = 376
$data[39,5] = 0
$data[39,6] = 0
$data[39,7] = 0
$data[40,0] = 'Eslovenia'
$data[40,1] = 383
$data[40,2] = 42
$data[40,3] = 0
$data[40,4] = 382
$data[40,5] = 12
$data[40,6] = 0
$data[40,7] = 1
$data[41,0] = 'Rumania'
$data[41,1] = 367
$data[41,2] = 59
$data[41,3] = 52
$data[41,4] = 315
$data[41,5] = 14
$data[41,6] = 0
$data[41,7] = 0
$data[42,0] = 'India'
$data[42,1] = 332
$data[42,2] = 83
$data[42,3] = 23
$data[42,4] = 304
$data[42,5] = 0
$data[42,6] = 0
$data[42,7] = 5
$data[43,0] = 'Peru'
$data[43,1] = 318
$data[43,2] = 55
$data[43,3] = 1
$data[43,4] = 312
$data[43,5] = 5
$data[43,6] = 1
$data[43,7] = 5
$data[44,0] = 'Barein'
$data[44,1] = 310
$data[44,2] = 12
$data[44,3] = 125
$data[44,4] = 184
$data[44,5] = 4
$data[44,6] = 0
$data[44,7] = 1
$data[45,0] = 'Filipinas'
$data[45,1] = 307
$data[45,2] = 77
$data[45,3] = 13
$data[45,4] = 275
$data[45,5] = 1
$data[45,6] = 1
$data[45,7] = 0
$data[46,0] = 'Estonia'
$data[46,1] = 306
$data[46,2] = 23
$data[46,3] = 2
$data[46,4] = 304
$data[46,5] = 1
$data[46,6] = 0
$data[46,7] = 0
$data[47,0] = 'Rusia'
$data[47,1] = 306
$data[47,2] = 53
$data[47,3] = 16
$data[47,4] = 289
$data[47,5] = 0
$data[47,6] = 0
$data[47,7] = 1
$data[48,0] = 'Egipto'
$data[48,1] = 294
$data[48,2] = 9
$data[48,3] = 42
$data[48,4] = 242
$data[48,5] = 0
$data[48,6] = 2
$data[48,7] = 10
$data[49,0] = 'Hong Kong'
$data[49,1] = 273
$data[49,2] = 17
$data[49,3] = 98
$data[49,4] = 171
$data[49,5] = 4
$data[49,6] = 0
$data[49,7] = 4
$data[50,0] = 'Sudafrica'
$data[50,1] = 240
$data[50,2] = 38
$data[50,3] = 2
$data[50,4] = 238
$data[50,5] = 0
$data[50,6] = 0
$data[50,7] = 0
$data[51,0] = 'Libano'
$data[51,1] = 230
$data[51,2] = 53
$data[51,3] = 8
$data[51,4] = 218
$data[51,5] = 4
$data[51,6] = 0
$data[51,7] = 4
$data[52,0] = 'Irak'
$data[52,1] = 214
$data[52,2] = 6
$data[52,3] = 51
$data[52,4] = 146
$data[52,5] = 0
$data[52,6] = 0
$data[52,7] = 17
$data[53,0] = 'Croacia'
$data[53,1] = 206
$data[53,2] = 76
$data[53,3] = 5
$data[53,4] = 200
$data[53,5] = 0
$data[53,6] = 0
$data[53,7] = 1
$data[54,0] = 'Mexico'
$data[54,1] = 203
$data[54,2] = 39
$data[54,3] = 4
$data[54,4] = 197
$data[54,5] = 1
$data[54,6] = 1
$data[54,7] = 2
$data[55,0] = 'Panama'
$data[55,1] = 200
$data[55,2] = 0
$data[55,3] = 1
$data[55,4] = 198
$data[55,5] = 7
$data[55,6] = 0
$data[55,7] = 1
$data[56,0] = 'Colombia'
$data[56,1] = 196
$data[56,2] = 51
$data[56,3] = 1
$data[56,4] = 195
$data[56,5] = 0
$data[56,6] = 0
$data[56,7] = 0
$data[57,0] = 'Eslovaquia'
$data[57,1] = 178
$data[57,2] = 41
$data[57,3] = 7
$data[57,4] = 171
$data[57,5] = 2
$data[57,6] = 0
$data[57,7] = 0
$data[58,0] = 'Kuwait'
$data[58,1] = 176
$data[58,2] = 17
$data[58,3] = 27
$data[58,4] = 149
$data[58,5] = 5
$data[58,6] = 0
$data[58,7] = 0
$data[59,0] = 'Serbia'
$data[59,1] = 171
$data[59,2] = 36
$data[59,3] = 2
$data[59,4] = 168
$data[59,5] = 4
$data[59,6] = 0
$data[59,7] = 1
$data[60,0] = 'Bulgaria'
$data[60,1] = 163
$data[60,2] = 36
$data[60,3] = 3
$data[60,4] = 157
$data[60,5] = 3
$data[60,6] = 0
$data[60,7] = 3
$data[61,0] = 'Armenia'
$data[61,1] = 160
$data[61,2] = 24
$data[61,3] = 1
$data[61,4] = 159
$data[61,5] = 2
$data[61,6] = 0
$data[61,7] = 0
$data[62,0] = 'Argentina'
$data[62,1] = 158
$data[62,2] = 0
$data[62,3] = 3
$data[62,4] = 151
$data[62,5] = 0
$data[62,6] = 1
$data[62,7] = 4
$data[63,0] = 'Taiwan'
$data[63,1] = 153
$data[63,2] = 18
$data[63,3] = 28
$data[63,4] = 123
$data[63,5] = 0
$data[63,6] = 0
$data[63,7] = 2
$data[64,0] = 'Emiratos Arabes Unidos'
$data[64,1] = 153
$data[64,2] = 13
$data[64,3] = 38
$data[64,4] = 113
$data[64,5] = 2
$data[64,6] = 0
$data[64,7] = 2
$data[65,0] = 'San Marino'
$data[65,1] = 151
$data[65,2] = 0
$data[65,3] = 4
$data[65,4] = 127
$data[65,5] = 12
$data[65,6] = 6
$data[65,7] = 20
$data[66,0] = 'Argelia'
$data[66,1] = 139
$data[66,2] = 45
$data[66,3] = 65
$data[66,4] = 59
$data[66,5] = 0
$data[66,6] = 4
$data[66,7] = 15
$data[67,0] = 'Letonia'
$data[67,1] = 124
$data[67,2] = 13
$data[67,3] = 1
$data[67,4] = 123
$data[67,5] = 0
$data[67,6] = 0
$data[67,7] = 0
$data[68,0] = 'Costa Rica'
$data[68,1] = 117
$data[68,2] = 4
$data[68,3] = 2
$data[68,4] = 113
$data[68,5] = 2
$data[68,6] = 0
$data[68,7] = 2
$data[69,0] = 'Republica Dominicana'
$data[69,1] = 112
$data[69,2] = 40
$data[69,3] = 0
$data[69,4] = 109
$data[69,5] = 0
$data[69,6] = 1
$data[69,7] = 3
$data[70,0] = 'Uruguay'
$data[70,1] = 110
$data[70,2] = 0
$data[70,3] = 0
$data[70,4] = 110
$data[70,5] = 0
$data[70,6] = 0
$data[70,7] = 0
$data[71,0] = 'Hungria'
$data[71,1] = 103
$data[71,2] = 18
$data[71,3] = 7
$data[71,4] = 92
$data[71,5] = 6
$data[71,6] = 0
$data[71,7] = 4
$data[72,0] = 'Jordania'
$data[72,1] = 99
$data[72,2] = 15
$data[72,3] = 1
$data[72,4] = 98
$data[72,5] = 0
$data[72,6] = 0
$data[72,7] = 0
$data[73,0] = 'Vietnam'
$data[73,1] = 94
$data[73,2] = 3
$data[73,3] = 17
$data[73,4] = 77
$data[73,5] = 2
$data[73,6] = 0
$data[73,7] = 0
$data[74,0] = 'Bosnia y Herzegovina'
$data[74,1] = 93
$data[74,2] = 4
$data[74,3] = 2
$data[74,4] = 90
$data[74,5] = 1
$data[74,6] = 1
$data[74,7] = 1
$data[75,0] = 'Islas Feroe'
$data[75,1] = 92
$data[75,2] = 12
$data[75,3] = 3
$data[75,4] = 89
$data[75,5] = 0
$data[75,6] = 0
$data[75,7] = 0
$data[76,0] = 'Principado de Andorra'
$data[76,1] = 88
$data[76,2] = 13
$data[76,3] = 1
$data[76,4] = 87
$data[76,5] = 2
$data[76,6] = 0
$data[76,7] = 0
$data[77,0] = 'Marruecos'
$data[77,1] = 86
$data[77,2] = 0
$data[77,3] = 2
$data[77,4] = 81
$data[77,5] = 1
$data[77,6] = 0
$data[77,7] = 3
$data[78,0] = 'Republica de Macedonia'
$data[78,1] = 85
$data[78,2] = 9
$data[78,3] = 1
$data[78,4] = 84
$data[78,5] = 1
$data[78,6] = 0
$data[78,7] = 0
$data[79,0] = 'Republica de Chipre'
$data[79,1] = 84
$data[79,2] = 9
$data[79,3] = 0
$data[79,4] = 84
$data[79,5] = 1
$data[79,6] = 0
$data[79,7] = 0
$data[80,0] = 'Brunei'
$data[80,1] = 83
$data[80,2] = 5
$data[80,3] = 1
$data[80,4] = 82
$data[80,5] = 2
$data[80,6] = 0
$data[80,7] = 0
$data[81,0] = 'Lituania'
$data[81,1] = 83
$data[81,2] = 14
$data[81,3] = 1
$data[81,4] = 81
$data[81,5] = 1
$data[81,6] = 0
$data[81,7] = 1
$data[82,0] = 'Moldavia'
$data[82,1] = 80
$data[82,2] = 14
$data[82,3] = 1
$data[82,4] = 78
$data[82,5] = 3
$data[82,6] = 0
$data[82,7] = 1
$data[83,0] = 'Sri Lanka'
$data[83,1] = 77
$data[83,2] = 4
$data[83,3] = 3
$data[83,4] = 74
$data[83,5] = 2
$data[83,6] = 0
$data[83,7] = 0
$data[84,0] = 'Albania'
$data[84,1] = 76
$data[84,2] = 6
$data[84,3] = 2
$data[84,4] = 72
$data[84,5] = 2
$data[84,6] = 0
$data[84,7] = 2
$data[85,0] = 'Bielorrusia'
$data[85,1] = 76
$data[85,2] = 7
$data[85,3] = 15
$data[85,4] = 61
$data[85,5] = 0
$data[85,6] = 0
$data[85,7] = 0
$data[86,0] = 'Malta'
$data[86,1] = 73
$data[86,2] = 9
$data[86,3] = 2
$data[86,4] = 71
$data[86,5] = 1
$data[86,6] = 0
$data[86,7] = 0
$data[87,0] = 'Venezuela'
$data[87,1] = 70
$data[87,2] = 5
$data[87,3] = 15
$data[87,4] = 55
$data[87,5] = 2
$data[87,6] = 0
$data[87,7] = 0
$data[88,0] = 'Burkina Faso'
$data[88,1] = 64
$data[88,2] = 24
$data[88,3] = 5
$data[88,4] = 56
$data[88,5] = 0
$data[88,6] = 2
$data[88,7] = 3
$data[89,0] = 'Tunez'
$data[89,1] = 60
$data[89,2] = 6
$data[89,3] = 1
$data[89,4] = 58
$data[89,5] = 7
$data[89,6] = 0
$data[89,7] = 1
$data[90,0] = 'Kazajistan'
$data[90,1] = 53
$data[90,2] = 1
$data[90,3] = 0
$data[90,4] = 53
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 0
$data[91,0] = 'Camboya'
$data[91,1] = 53
$data[91,2] = 2
$data[91,3] = 2
$data[91,4] = 51
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 0
$data[92,0] = 'Azerbaiyan'
$data[92,1] = 53
$data[92,2] = 9
$data[92,3] = 11
$data[92,4] = 41
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 1
$data[93,0] = 'Nueva Zelanda'
$data[93,1] = 52
$data[93,2] = 0
$data[93,3] = 0
$data[93,4] = 52
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 0
$data[94,0] = 'Oman'
$data[94,1] = 52
$data[94,2] = 4
$data[94,3] = 13
$data[94,4] = 39
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 0
$data[95,0] = 'Estado de Palestina'
$data[95,1] = 52
$data[95,2] = 4
$data[95,3] = 17
$data[95,4] = 35
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 0
$data[96,0] = 'Guadalupe'
$data[96,1] = 51
$data[96,2] = 0
$data[96,3] = 0
$data[96,4] = 50
$data[96,5] = 4
$data[96,6] = 0
$data[96,7] = 1
$data[97,0] = 'Trinidad yTobago'
$data[97,1] = 49
$data[97,2] = 40
$data[97,3] = 0
$data[97,4] = 49
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 0
$data[98,0] = 'Georgia'
$data[98,1] = 49
$data[98,2] = 5
$data[98,3] = 1
$data[98,4] = 48
$data[98,5] = 1
$data[98,6] = 0
$data[98,7] = 0
$data[99,0] = 'Ucrania'
$data[99,1] = 47
$data[99,2] = 6
$data[99,3] = 1
$data[99,4] = 43
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 3
$data[100,0] = 'Senegal'
$data[100,1] = 47
$data[100,2] = 0
$data[100,3] = 5
$data[100,4] = 42
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 0
$data[101,0] = 'Reunion'
$data[101,1] = 45
$data[101,2] = 7
$data[101,3] = 0
$data[101,4] = 45
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 0
$data[102,0] = 'Uzbekistan'
$data[102,1] = 41
$data[102,2] = 8
$data[102,3] = 0
$data[102,4] = 41
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 0
$data[103,0] = 'Liechtenstein'
$data[103,1] = 37
$data[103,2] = 9
$data[103,3] = 0
$data[103,4] = 37
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 0
$data[104,0] = 'Martinica'
$data[104,1] = 32
$data[104,2] = 0
$data[104,3] = 0
$data[104,4] = 31
$data[104,5] = 7
$data[104,6] = 0
$data[104,7] = 1
$data[105,0] = 'Camerun'
$data[105,1] = 27
$data[105,2] = 0
$data[105,3] = 2
$data[105,4] = 25
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 0
$data[106,0] = 'Honduras'
$data[106,1] = 24
$data[106,2] = 0
$data[106,3] = 0
$data[106,4] = 24
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 0
$data[107,0] = 'Afganistan'
$data[107,1] = 24
$data[107,2] = 0
$data[107,3] = 1
$data[107,4] = 23
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 0
$data[108,0] = 'Banglades'
$data[108,1] = 24
$data[108,2] = 4
$data[108,3] = 3
$data[108,4] = 19
$data[108,5] = 0
$data[108,6] = 1
$data[108,7] = 2
$data[109,0] = 'Consejo Danes para los Refugiados'
$data[109,1] = 23
$data[109,2] = 5
$data[109,3] = 0
$data[109,4] = 23
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 0
$data[110,0] = 'Nigeria'
$data[110,1] = 22
$data[110,2] = 10
$data[110,3] = 1
$data[110,4] = 21
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 0
$data[111,0] = 'Cuba'
$data[111,1] = 21
$data[111,2] = 0
$data[111,3] = 0
$data[111,4] = 20
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 1
$data[112,0] = 'Ghana'
$data[112,1] = 19
$data[112,2] = 3
$data[112,3] = 0
$data[112,4] = 19
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 0
$data[113,0] = 'Bolivia'
$data[113,1] = 19
$data[113,2] = 3
$data[113,3] = 0
$data[113,4] = 19
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 0
$data[114,0] = 'Jamaica'
$data[114,1] = 19
$data[114,2] = 0
$data[114,3] = 2
$data[114,4] = 16
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 1
$data[115,0] = 'Guayana Francesa'
$data[115,1] = 18
$data[115,2] = 3
$data[115,3] = 0
$data[115,4] = 18
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 0
$data[116,0] = 'Paraguay'
$data[116,1] = 18
$data[116,2] = 0
$data[116,3] = 0
$data[116,4] = 17
$data[116,5] = 1
$data[116,6] = 1
$data[116,7] = 1
$data[117,0] = 'Macao'
$data[117,1] = 18
$data[117,2] = 1
$data[117,3] = 10
$data[117,4] = 8
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 0
$data[118,0] = 'Puerto Rico'
$data[118,1] = 17
$data[118,2] = 3
$data[118,3] = 0
$data[118,4] = 17
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 0
$data[119,0] = 'Ruanda'
$data[119,1] = 17
$data[119,2] = 0
$data[119,3] = 0
$data[119,4] = 17
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 0
$data[120,0] = 'Guatemala'
$data[120,1] = 17
$data[120,2] = 5
$data[120,3] = 0
$data[120,4] = 16
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 1
$data[121,0] = 'Togo'
$data[121,1] = 16
$data[121,2] = 7
$data[121,3] = 0
$data[121,4] = 16
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 0
$data[122,0] = 'Guam'
$data[122,1] = 15
$data[122,2] = 1
$data[122,3] = 0
$data[122,4] = 15
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 0
$data[123,0] = 'Polinesia Francesa'
$data[123,1] = 15
$data[123,2] = 4
$data[123,3] = 0
$data[123,4] = 15
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = 'Kirguistan'
$data[124,1] = 14
$data[124,2] = 8
$data[124,3] = 0
$data[124,4] = 14
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = 'Montenegro'
$data[125,1] = 14
$data[125,2] = 0
$data[125,3] = 0
$data[125,4] = 14
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = 'Costa de Marfil'
$data[126,1] = 14
$data[126,2] = 5
$data[126,3] = 1
$data[126,4] = 13
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 0
$data[127,0] = 'Mauricio'
$data[127,1] = 14
$data[127,2] = 2
$data[127,3] = 0
$data[127,4] = 13
$data[127,5] = 0
$data[127,6] = 1
$data[127,7] = 1
$data[128,0] = 'Maldivas'
$data[128,1] = 13
$data[128,2] = 0
$data[128,3] = 3
$data[128,4] = 10
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = 'Monaco'
$data[129,1] = 11
$data[129,2] = 0
$data[129,3] = 0
$data[129,4] = 11
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = 'Mongolia'
$data[130,1] = 10
$data[130,2] = 4
$data[130,3] = 0
$data[130,4] = 10
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = 'Gibraltar'
$data[131,1] = 10
$data[131,2] = 0
$data[131,3] = 2
$data[131,4] = 8
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 0
$data[132,0] = 'Etiopia'
$data[132,1] = 9
$data[132,2] = 0
$data[132,3] = 0
$data[132,4] = 9
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = 'Kenia'
$data[133,1] = 7
$data[133,2] = 0
$data[133,3] = 0
$data[133,4] = 7
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = 'Seychelles'
$data[134,1] = 7
$data[134,2] = 0
$data[134,3] = 0
$data[134,4] = 7
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 0
$data[135,0] = 'Mayotte'
$data[135,1] = 7
$data[135,2] = 0
$data[135,3] = 0
$data[135,4] = 7
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = 'Tanzania'
$data[136,1] = 6
$data[136,2] = 0
$data[136,3] = 0
$data[136,4] = 6
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = 'Guinea Ecuatorial'
$data[137,1] = 6
$data[137,2] = 0
$data[137,3] = 0
$data[137,4] = 6
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 0
$data[138,0] = 'Barbados'
$data[138,1] = 6
$data[138,2] = 0
$data[138,3] = 0
$data[138,4] = 6
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = 'Islas Virgenes de los Estados Unidos'
$data[139,1] = 6
$data[139,2] = 3
$data[139,3] = 0
$data[139,4] = 6
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = 'Guyana'
$data[140,1] = 5
$data[140,2] = 0
$data[140,3] = 0
$data[140,4] = 4
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 1
$data[141,0] = 'Aruba'
$data[141,1] = 5
$data[141,2] = 0
$data[141,3] = 1
$data[141,4] = 4
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = 'Nueva Caledonia'
$data[142,1] = 4
$data[142,2] = 2
$data[142,3] = 0
$data[142,4] = 4
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = 'Surinam'
$data[143,1] = 4
$data[143,2] = 0
$data[143,3] = 0
$data[143,4] = 4
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = 'San Martin (Parte Francesa)'
$data[144,1] = 4
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 4
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = 'Bahamas'
$data[145,1] = 4
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 4
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 0
$data[146,0] = 'Gabon'
$data[146,1] = 4
$data[146,2] = 0
$data[146,3] = 0
$data[146,4] = 3
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 1
$data[147,0] = 'El Salvador'
$data[147,1] = 3
$data[147,2] = 2
$data[147,3] = 0
$data[147,4] = 3
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 0
$data[148,0] = 'Republica de Africa Central'
$data[148,1] = 3
$data[148,2] = 0
$data[148,3] = 0
$data[148,4] = 3
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 0
$data[149,0] = 'Zimbabue'
$data[149,1] = 3
$data[149,2] = 2
$data[149,3] = 0
$data[149,4] = 3
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 0
$data[150,0] = 'San Bartolome'
$data[150,1] = 3
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 3
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 0
$data[151,0] = 'Liberia'
$data[151,1] = 3
$data[151,2] = 1
$data[151,3] = 0
$data[151,4] = 3
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 0
$data[152,0] = 'Congo'
$data[152,1] = 3
$data[152,2] = 0
$data[152,3] = 0
$data[152,4] = 3
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0
$data[153,0] = 'Namibia'
$data[153,1] = 3
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 3
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = 'Cabo Verde'
$data[154,1] = 3
$data[154,2] = 2
$data[154,3] = 0
$data[154,4] = 3
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = 'Madagascar'
$data[155,1] = 3
$data[155,2] = 0
$data[155,3] = 0
$data[155,4] = 3
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = 'Curazao'
$data[156,1] = 3
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 2
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 1
$data[157,0] = 'Islas Caimanes'
$data[157,1] = 3
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 2
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 1
$data[158,0] = 'Isla de Man'
$data[158,1] = 2
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 2
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = 'Bermudas'
$data[159,1] = 2
$data[159,2] = 0
$data[159,3] = 0
$data[159,4] = 2
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0
$data[160,0] = 'Santa Lucia'
$data[160,1] = 2
$data[160,2] = 0
$data[160,3] = 0
$data[160,4] = 2
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = 'Fiyi'
$data[161,1] = 2
$data[161,2] = 1
$data[161,3] = 0
$data[161,4] = 2
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = 'Nicaragua'
$data[162,1] = 2
$data[162,2] = 0
$data[162,3] = 0
$data[162,4] = 2
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = 'Butan'
$data[163,1] = 2
$data[163,2] = 0
$data[163,3] = 0
$data[163,4] = 2
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = 'Zambia'
$data[164,1] = 2
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 2
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = 'Haiti'
$data[165,1] = 2
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 2
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 0
$data[166,0] = 'Groenlandia'
$data[166,1] = 2
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 2
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = 'Benin'
$data[167,1] = 2
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 2
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = 'Angola'
$data[168,1] = 2
$data[168,2] = 1
$data[168,3] = 0
$data[168,4] = 2
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = 'Mauritania'
$data[169,1] = 2
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 2
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = 'Guinea'
$data[170,1] = 2
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 2
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = 'Sudan'
$data[171,1] = 2
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 1
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 1
$data[172,0] = 'Suazilandia'
$data[172,1] = 1
$data[172,2] = 0
$data[172,3] = 0
$data[172,4] = 1
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = 'Eritrea'
$data[173,1] = 1
$data[173,2] = 1
$data[173,3] = 0
$data[173,4] = 1
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = 'Santa Sede'
$data[174,1] = 1
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 1
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = 'San Martin (Parte Holandesa)'
$data[175,1] = 1
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 1
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = 'Republica del Chad'
$data[176,1] = 1
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 1
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = 'Niger'
$data[177,1] = 1
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 1
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = 'Montserrat'
$data[178,1] = 1
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 1
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = 'Gambia'
$data[179,1] = 1
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 1
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = 'Antigua y Barbuda'
$data[180,1] = 1
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 1
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = 'San Vicente y las Granadinas'
$data[181,1] = 1
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 1
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = 'Republica de Yibuti'
$data[182,1] = 1
$data[182,2] = 0
$data[182,3] = 0
$data[182,4] = 1
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Papua Nueva Guinea'
$data[183,1] = 1
$data[183,2] = 0
$data[183,3] = 0
$data[183,4] = 1
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = 'Timor Oriental'
$data[184,1] = 1
$data[184,2] = 1
$data[184,3] = 0
$data[184,4] = 1
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'Somalia'
$data[185,1] = 1
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 1
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = 'Nepal'
$data[186,1] = 1
$data[186,2] = 0
$data[186,3] = 1
$data[186,4] = 0
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0

$ws.Range("A4:H190").Value = $data
